$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1692045.6
$ws.Range("J17").Value = 1728793.1
$ws.Range("L17").Value = 5186379.300000001
$ws.Range("N17").Value = -5186715.300000001
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H81").Value = 75995
$ws.Range("J81").Value = 75995
$ws.Range("L81").Value = 75995
$ws.Range("N81").Value = -77991
$ws.Range("H84").Value = 75995
$ws.Range("J84").Value = 75995
$ws.Range("L84").Value = 227985
$ws.Range("N84").Value = -237969
$ws.Range("H86").Value = 5110.8
$ws.Range("I86").Value = 1529.2
$ws.Range("K86").Value = 1529.2
$ws.Range("M86").Value = -406.2
$ws.Range("H89").Value = 5110.8
$ws.Range("I89").Value = 1529.2
$ws.Range("K89").Value = 7646
$ws.Range("M89").Value = -2030
$ws.Range("H107").Value = 673.6667
$ws.Range("J107").Value = 195.33333
$ws.Range("L107").Value = 195.33333
$ws.Range("N107").Value = -4035.33333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 88560.71000000001
$ws.Range("J37").Value = 88560.71000000001
$ws.Range("L37").Value = 88560.71000000001
$ws.Range("N37").Value = -89106.71000000001
$ws.Range("H74").Value = 2221
$ws.Range("I74").Value = 1829.12
$ws.Range("K74").Value = 1829.12
$ws.Range("M74").Value = -955.1199999999999
$ws.Range("H77").Value = 2221
$ws.Range("I77").Value = 1829.12
$ws.Range("K77").Value = 9145.599999999999
$ws.Range("M77").Value = -4777.599999999999
$ws.Range("H88").Value = 1741.3158
$ws.Range("J88").Value = 1768.5714
$ws.Range("L88").Value = 1768.5714
$ws.Range("N88").Value = -2580.5714
$ws.Range("H91").Value = 1741.3158
$ws.Range("J91").Value = 1768.5714
$ws.Range("L91").Value = 1768.5714
$ws.Range("N91").Value = -4576.5714
$ws.Range("H97").Value = 1514
$ws.Range("I97").Value = 1514
$ws.Range("K97").Value = 1514
$ws.Range("M97").Value = -1018
$ws.Range("H122").Value = 3071.75
$ws.Range("I122").Value = 2680.64
$ws.Range("K122").Value = 8041.92
$ws.Range("M122").Value = -5591.92
$ws.Range("H132").Value = 1899.9474
$ws.Range("I132").Value = 1495.3823
$ws.Range("K132").Value = 4486.1469
$ws.Range("M132").Value = -1956.1469

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 879.1111
$ws.Range("I94").Value = 201.71428
$ws.Range("J94").Value = 3250
$ws.Range("K94").Value = 201.71428
$ws.Range("L94").Value = 3250
$ws.Range("M94").Value = 249.28572
$ws.Range("N94").Value = -4152
$ws.Range("H105").Value = 3541.6562
$ws.Range("I105").Value = 2348.963
$ws.Range("K105").Value = 2348.963
$ws.Range("M105").Value = -601.9630000000002
$ws.Range("H134").Value = 7200.567
$ws.Range("I134").Value = 2365.182
$ws.Range("K134").Value = 7095.545999999999
$ws.Range("M134").Value = -4560.545999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 566.1177
$ws.Range("I22").Value = 538.36365
$ws.Range("K22").Value = 538.36365
$ws.Range("M22").Value = -188.36365
$ws.Range("H99").Value = 2565.9524
$ws.Range("J99").Value = 3499.6667
$ws.Range("L99").Value = 3499.6667
$ws.Range("N99").Value = -6495.6667
$ws.Range("H126").Value = 2565.9524
$ws.Range("J126").Value = 3499.6667
$ws.Range("L126").Value = 10499.0001
$ws.Range("N126").Value = -15439.0001
$ws.Range("H134").Value = 2792.85
$ws.Range("J134").Value = 5807.4614
$ws.Range("L134").Value = 17422.3842
$ws.Range("N134").Value = -22492.3842

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 862.3333
$ws.Range("I117").Value = 551.2
$ws.Range("J117").Value = 1484.6
$ws.Range("K117").Value = 1653.6
$ws.Range("L117").Value = 4453.799999999999
$ws.Range("M117").Value = 1788.4
$ws.Range("N117").Value = -11337.8
$ws.Range("H121").Value = 4210.1
$ws.Range("I121").Value = 483
$ws.Range("J121").Value = 5807.4287
$ws.Range("K121").Value = 1449
$ws.Range("L121").Value = 17422.2861
$ws.Range("M121").Value = -139
$ws.Range("N121").Value = -20042.2861
$ws.Range("H123").Value = 4513
$ws.Range("I123").Value = 3673
$ws.Range("J123").Value = 7033
$ws.Range("K123").Value = 11019
$ws.Range("L123").Value = 21099
$ws.Range("M123").Value = -8569
$ws.Range("N123").Value = -25999
$ws.Range("H125").Value = 3537.2856
$ws.Range("I125").Value = 2632.6
$ws.Range("K125").Value = 7897.799999999999
$ws.Range("M125").Value = -2977.799999999999
$ws.Range("H126").Value = 7885.25
$ws.Range("J126").Value = 9847
$ws.Range("L126").Value = 29541
$ws.Range("N126").Value = -39421
$ws.Range("H131").Value = 1749
$ws.Range("I131").Value = 965
$ws.Range("J131").Value = 1823.6666
$ws.Range("K131").Value = 2895
$ws.Range("L131").Value = 5470.9998
$ws.Range("M131").Value = 2145
$ws.Range("N131").Value = -15550.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2825.2856
$ws.Range("I97").Value = 1895.6
$ws.Range("J97").Value = 5149.5
$ws.Range("K97").Value = 1895.6
$ws.Range("L97").Value = 5149.5
$ws.Range("M97").Value = -1399.6
$ws.Range("N97").Value = -6141.5
$ws.Range("H102").Value = 30301.143
$ws.Range("I102").Value = 1381.8462
$ws.Range("K102").Value = 1381.8462
$ws.Range("M102").Value = 240.1538
$ws.Range("H105").Value = 500025000
$ws.Range("J105").Value = 500025000
$ws.Range("L105").Value = 500025000
$ws.Range("N105").Value = -500031988
$ws.Range("H122").Value = 2543.65
$ws.Range("I122").Value = 2284.8667
$ws.Range("J122").Value = 3320
$ws.Range("K122").Value = 6854.6001
$ws.Range("L122").Value = 9960
$ws.Range("M122").Value = -4404.6001
$ws.Range("N122").Value = -14860
$ws.Range("H131").Value = 76325
$ws.Range("J131").Value = 76325
$ws.Range("L131").Value = 76325
$ws.Range("N131").Value = -86405
$ws.Range("H132").Value = 12348931
$ws.Range("I132").Value = 13892016
$ws.Range("J132").Value = 4252.6665
$ws.Range("K132").Value = 41676048
$ws.Range("L132").Value = 12757.9995
$ws.Range("M132").Value = -41673518
$ws.Range("N132").Value = -17817.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9096619
$ws.Range("I7").Value = 20004212
$ws.Range("J7").Value = 6958
$ws.Range("K7").Value = 20004212
$ws.Range("L7").Value = 6958
$ws.Range("M7").Value = -20004100
$ws.Range("N7").Value = -7182
$ws.Range("H22").Value = 2871.3157
$ws.Range("I22").Value = 1286.4
$ws.Range("K22").Value = 1286.4
$ws.Range("M22").Value = -991.4000000000001
$ws.Range("H27").Value = 2871.3157
$ws.Range("I27").Value = 1286.4
$ws.Range("K27").Value = 1286.4
$ws.Range("M27").Value = -1179.4
$ws.Range("H95").Value = 52000
$ws.Range("J95").Value = 52000
$ws.Range("L95").Value = 52000
$ws.Range("N95").Value = -57492
$ws.Range("H115").Value = 89884.5
$ws.Range("J115").Value = 89884.5
$ws.Range("L115").Value = 89884.5
$ws.Range("N115").Value = -92234.5
$ws.Range("H126").Value = 9096619
$ws.Range("I126").Value = 20004212
$ws.Range("J126").Value = 6958
$ws.Range("K126").Value = 60012636
$ws.Range("L126").Value = 20874
$ws.Range("M126").Value = -60010166
$ws.Range("N126").Value = -25814
$ws.Range("H136").Value = 4135.815
$ws.Range("I136").Value = 3871.5908
$ws.Range("K136").Value = 11614.7724
$ws.Range("M136").Value = -9064.7724

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5000
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H100").Value = 2429.0557
$ws.Range("I100").Value = 1948.3334
$ws.Range("K100").Value = 3896.6668
$ws.Range("M100").Value = -3355.6668
$ws.Range("H132").Value = 2096.5745
$ws.Range("I132").Value = 2044.0541
$ws.Range("J132").Value = 2290.9
$ws.Range("K132").Value = 6132.1623
$ws.Range("L132").Value = 6872.700000000001
$ws.Range("M132").Value = -3602.1623
$ws.Range("N132").Value = -11932.7
$ws.Range("H136").Value = 14093
$ws.Range("I136").Value = 24245.445
$ws.Range("K136").Value = 72736.33499999999
$ws.Range("M136").Value = -70186.33499999999
